$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data range (A1:D4) entirely before rewriting
$ws.Range("A1:D4").Clear()

# Header row
$ws.Range("A1").Value = "Title "
$ws.Range("B1").Value = "Studio "
$ws.Range("C1").Value = "Likes"

# Data rows
# Note: value-setting order controls shared-string table allocation order,
# so "sa" (B6) must be written before "a" (B3) is first introduced.
$ws.Range("A2").Value = "dsad"
$ws.Range("B2").Value = "d"
$ws.Range("C2").Value = "dsad"

$ws.Range("B6").Value = "sa"

$ws.Range("B3").Value = "a"
$ws.Range("C3").Value = "dsad"

$ws.Range("B4").Value = "d"

$ws.Range("B5").Value = "a"

$ws.Range("G11").Select()
